# "Change calcul to <arrondie> in situation de cloture"
# Rewrites the data rows of the "Etat Taxes" sheet (taxe/loyer situation de
# cloture) with the recalculated ("arrondie") figures, and drops the two
# trailing summary/filler rows that are no longer part of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -----------------------------------------------------------
$ws.Range("A2").Value = "949/DR"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "BG12456"
$ws.Range("D2").Value = "HAYLALA ONE"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 4334.4
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 433.44
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 3900.96

# --- Row 3 -----------------------------------------------------------
$ws.Range("A3").Value = "949/DR"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "BG196435"
$ws.Range("D3").Value = "HAYLAL TWO"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 5665.73
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 566.57
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 5099.16

# --- Row 4 -------------------------------------------------------------
$ws.Range("A4").Value = "001/LF/TEST DR/AV1"
$ws.Range("B4").Value = "Logement de fonction"
# C4 is a purely-numeric string in the source data ("11986345") - force it
# to stay text (same logical type as the rest of column C) instead of
# letting it be auto-detected as a number, then drop back to the default
# style so no stray number-format sticks around on the cell.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "11986345"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "ALI EXPRESSE"
$ws.Range("E4").Value = "oui"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 30000
$ws.Range("I4").Value = 30000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = "--"
$ws.Range("M4").Value = 30000

# --- Row 5 -----------------------------------------------------------
$ws.Range("A5").Value = "001/TEST DR"
$ws.Range("B5").Value = "Direction régionale"
$ws.Range("C5").Value = "BG432432"
$ws.Range("D5").Value = "TETS TESTS"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 20000
$ws.Range("I5").Value = 20000
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = "--"
$ws.Range("M5").Value = 17000

# --- Row 6 (was the totals row on row 8) ------------------------------
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = 60000.13
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 4000.01
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 56000.12

# --- Drop the old rows 7 and 8 (data now fits in 6 rows) --------------
$ws.Rows("7:8").Delete() | Out-Null
